$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.556.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.87%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.049.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.11%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.12%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'384.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.40%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'102.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.84%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.544"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.06%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D9").Value = "'0.585"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.92%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'36.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.63%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.26%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.0863"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.42%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.540.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.61%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'18.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.36%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'7.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.85%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.075.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +4.26%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.972"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.12%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'10.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -5.75%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'51.631.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.96%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'3.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.80%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'12.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.01%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.0₃0964"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.44%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'70.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.03%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'268.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.56%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'3.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.42%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'8.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +7.56%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'27.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +4.60%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +5.26%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.01%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.09%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.107"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.32%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'10.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.05%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'34.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.86%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.34%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'50.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.29%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.0444"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.42%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.01%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'3.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +3.87%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.290"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +7.86%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'16.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.35%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.38%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'128.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.93%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.41%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'2.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.96%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'3.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +4.27%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'21.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.50%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +5.29%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +2.79%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.042.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.50%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'3.361.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.37%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +7.37%  "
$ws.Range("E51").Style = "Normal"
